$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per corrected consignas
$ws.Range("F2").Value = " F,"
$ws.Range("I2").Value = " Xy,"
$ws.Range("J2").Value = "?"
$ws.Range("M2").Value = "una radiografia"
$ws.Range("N2").Value = "?"

# Remove rows 3 and 4 entirely
$ws.Rows("3:4").Delete()
